# Update "want to go" counts (column F) that were refreshed from the source
# site in the regenerated gh-pages output.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 625
$wsExhibit.Range("F6").Value = 119
$wsExhibit.Range("F10").Value = 5014
$wsExhibit.Range("F11").Value = 4701
$wsExhibit.Range("F13").Value = 30
$wsExhibit.Range("F15").Value = 42
$wsExhibit.Range("F16").Value = 174

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 72

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 625
$wsAll.Range("F6").Value = 119
$wsAll.Range("F10").Value = 5014
$wsAll.Range("F11").Value = 4701
$wsAll.Range("F13").Value = 30
$wsAll.Range("F15").Value = 42
$wsAll.Range("F16").Value = 174
$wsAll.Range("F17").Value = 72
